$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attributes")

# Insert a new column B (becomes region2) and push pset_pn/cset_cn/attribute/value right.
$ws.Range("B:B").Insert()

# Insert a new column E (becomes commodity2), shifting attribute/value right.
$ws.Range("E:E").Insert()

# Row 2 headers
$ws.Range("B2").Value = "region2"
$ws.Range("D2").Value = "commodity"
$ws.Range("E2").Value = "commodity2"

# Row 3 (originally NORTH row) - add linked region SOUTH and commodity2 ELC
$ws.Range("B3").Value = "SOUTH"
$ws.Range("E3").Value = "ELC"

# Row 4 (originally SOUTH row) - add linked region NORTH and commodity2 ELC
$ws.Range("B4").Value = "NORTH"
$ws.Range("E4").Value = "ELC"

$wb.Save()
